$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear D3 value (was a numeric value, becomes an empty cell)
$ws.Range("D3").ClearContents()

# Row 7: rename "Other" to "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 224.7342885270675

# Copy style from A7 to A8 (bold + border + alignment) before filling A8
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

# New row 8: "Other" with a new value
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 99.13037193941483
